$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11 (shifts old row 11+ down by one)
$ws.Rows.Item(11).Insert()

# Copy formatting from row 10 into the new row 11 (font/border/numberformat etc.)
$ws.Range("A10:F10").Copy()
$ws.Range("A11:F11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set the row height to match the new content block
$ws.Rows.Item(11).RowHeight = 87.75

$A11 = "좌측 운전 차량의 운전자가 사용하는 디스플레이 및 컨트롤"

$D11 = @"
운전석 근처의 디스플레이 및 컨트롤 위치를 보여줍니다.
<br><h3>스티어링 휠 및 계기 패널</h3>
<br><img src="https://www.volvocars.com/images/support/img33f5a83a7a1f0d23c0a80152723389d1_1_--_--_VOICEpnghigh.png" width="422" height="291">
<br><img src="https://www.volvocars.com/images/support/img0fdb84773e1be862c0a801526d85d772_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
차폭등, 주간 주행등, 하향등, 상향등, 방향지시등, 전방 안개등/코너링 라이트*, 후방 안개등, 구간거리계 재설정
<br><img src="https://www.volvocars.com/images/support/imgd145437d66762d90c0a801520501e45a_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
헤드업 디스플레이
<br><img src="https://www.volvocars.com/images/support/img4fb3580a66931339c0a801523fc8b2d6_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
운전자 화면
<br><img src="https://www.volvocars.com/images/support/img2e171f2d66927c84c0a80152217536c8_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
와이퍼 및 워셔, 레인 센서
<br><img src="https://www.volvocars.com/images/support/img13271b6666920bd6c0a801525a3d0845_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
스티어링 휠 우측 키패드
<br><img src="https://www.volvocars.com/images/support/img2d02abd46691bf4cc0a801524ebdde35_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
스티어링 휠 조절
<br><img src="https://www.volvocars.com/images/support/imgbcad23bb669173cbc0a80152453f9f4a_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
경적
<br><img src="https://www.volvocars.com/images/support/imgce0bdadb669126c3c0a8015208ab2f89_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
스티어링 휠 좌측 키패드
<br><img src="https://www.volvocars.com/images/support/img31249d6e6690ce66c0a8015207bd6db2_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
보닛 열기
<br><img src="https://www.volvocars.com/images/support/img95520c73669080d3c0a8015244b76916_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
테일게이트 잠금 해제/열기/닫기
<br><h3>루프 콘솔</h3>
<br><img src="https://www.volvocars.com/images/support/img6bc2c634cf23b937c0a801520cf42fb4_1_--_--_VOICEpnghigh.png" width="422" height="291">
<br><img src="https://www.volvocars.com/images/support/img0fdb84773e1be862c0a801526d85d772_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
파노라마 선루프
<br><img src="https://www.volvocars.com/images/support/imgd145437d66762d90c0a801520501e45a_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
앞좌석 독서등 및 실내 조명
<br><img src="https://www.volvocars.com/images/support/img4fb3580a66931339c0a801523fc8b2d6_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
운전자 화면
<br><img src="https://www.volvocars.com/images/support/img2e171f2d66927c84c0a80152217536c8_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
루프 콘솔의 화면, <img src="https://www.volvocars.com/images/support/imgc3aacf00dfc259b5c0a801525ce4fb32_1_--_--_VOICEpnghigh.png" width="19" height="19">버튼
<br><img src="https://www.volvocars.com/images/support/img13271b6666920bd6c0a801525a3d0845_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
실내 미러의 수동 디밍
<br><h3>중앙 및 터널 콘솔</h3>
<br><img src="https://www.volvocars.com/images/support/img9b5437ea3fd57ad6c0a80152416d240d_1_--_--_VOICEpnghigh.png" width="422" height="291">
<br><img src="https://www.volvocars.com/images/support/img0fdb84773e1be862c0a801526d85d772_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
중앙 화면
<br><img src="https://www.volvocars.com/images/support/imgd145437d66762d90c0a801520501e45a_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
비상등, 서리 제거, 미디어, 글로브 박스 열기
<br><img src="https://www.volvocars.com/images/support/img4fb3580a66931339c0a801523fc8b2d6_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
기어 선택 레버
<br><img src="https://www.volvocars.com/images/support/img2e171f2d66927c84c0a80152217536c8_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
시동 버튼
<br><img src="https://www.volvocars.com/images/support/img13271b6666920bd6c0a801525a3d0845_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
주차 브레이크
<br><img src="https://www.volvocars.com/images/support/img2d02abd46691bf4cc0a801524ebdde35_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
정지 시 자동 제동
<br><h3>운전석 도어</h3>
<br><img src="https://www.volvocars.com/images/support/imgc2f4760b7a28408ac0a80152163b9504_1_--_--_VOICEpnghigh.png" width="422" height="291">
<br><img src="https://www.volvocars.com/images/support/img0fdb84773e1be862c0a801526d85d772_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
전동 앞좌석용 메모리, 도어 미러 및 헤드업 디스플레이 설정
<br><img src="https://www.volvocars.com/images/support/imgd145437d66762d90c0a801520501e45a_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
중앙 잠금장치
<br><img src="https://www.volvocars.com/images/support/img4fb3580a66931339c0a801523fc8b2d6_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
전동 윈도우, 도어 미러 및 어린이 보호용 잠금장치
<br><img src="https://www.volvocars.com/images/support/img2e171f2d66927c84c0a80152217536c8_1_--_--_VOICEpnghigh.png" width="19" height="19"> 
앞좌석 조절 스위치
"@

$ws.Range("A11").Value = $A11
$ws.Range("B11").Value = 45463
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = $D11
$ws.Range("E11").Value = "All"
$ws.Range("F11").Value = ""

# Update sheet view to match the new selection/scroll position
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("A11:D11").Select()
